# Weekly driver report update for 2025-04-19
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Bad Drivers" table ---
# Row 3: Intel(R) Wi-Fi 6 AX201 160MHz - 23.90.0.2 (Critical Minutes + Good Roaming % change)
$ws.Range("C3").Value = 2920
$ws.Range("D3").Value = 88.2

# Row 4: MediaTek MT7921 Wi-Fi 6 802.11ax PCIe Adapter - 3.0.1.1297 (Critical Minutes + Good Roaming % change)
$ws.Range("C4").Value = 219
$ws.Range("D4").Value = 97.7

# Row 5: Totals (Critical Minutes changes)
$ws.Range("C5").Value = 3139

# --- "Good Drivers" table (rows 13-20) ---
# Row 13: now MediaTek 3.0.1.1255, samples 23159, 99.9% (unchanged), vintage cleared
$ws.Range("A13").Value = "MediaTek MT7921 Wi-Fi 6 802.11ax PCIe Adapter - 3.0.1.1255"
$ws.Range("B13").Value = 23159
$ws.Range("E13").Value = $null

# Row 14: now Intel 21.60.2.1, samples 56018, 100%, vintage cleared
$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B14").Value = 56018
$ws.Range("D14").Value = 100
$ws.Range("E14").Value = $null

# Row 15: now Intel 22.50.1.1, samples 34244, 100% (unchanged), vintage cleared
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B15").Value = 34244
$ws.Range("E15").Value = $null

# Row 16: now Intel 23.100.0.4, samples 442178, 99.9% (unchanged), vintage 2024-11-10
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B16").Value = 442178
$ws.Range("E16").Formula = '="2024-11-10"'

# Row 17: now Intel 22.80.0.9, samples 77849, 99.9%, vintage 2021-08-18
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B17").Value = 77849
$ws.Range("D17").Value = 99.90000000000001
$ws.Range("E17").Formula = '="2021-08-18"'

# Row 18: unchanged (Intel 21.110.3.2, samples 59673, 100%, vintage 2020-08-05)

# Row 19: now MediaTek 3.0.1.1216, samples 36106, 100%, vintage 2020-08-05
$ws.Range("A19").Value = "MediaTek MT7921 Wi-Fi 6 802.11ax PCIe Adapter - 3.0.1.1216"
$ws.Range("B19").Value = 36106
$ws.Range("E19").Formula = '="2020-08-05"'

# Row 20: now Intel 21.70.0.6, samples 113652, 100%, vintage 2019-12-14 (unchanged)
$ws.Range("A20").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B20").Value = 113652
